$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; this shifts the existing rows 71-146
# down to 72-147 (Excel also carries the row's number formatting, e.g.
# the date style on column D, down with them).
$ws.Rows("71:71").Insert()

# Populate the newly-inserted row 71 with the new weekly record.
$ws.Range("A71").Value = 6
$ws.Range("B71").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").Value = 44579
$ws.Range("E71").Value = 13
$ws.Range("F71").Value = 100112029
$ws.Range("G71").Value = "Orégano"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 34
$ws.Range("K71").Value = 8000
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = 8441
$ws.Range("N71").Value = "`$/docena de atados"
$ws.Range("O71").Value = "Región Metropolitana"
$ws.Range("P71").Value = 2814
$ws.Range("Q71").Value = 3
$ws.Range("R71").Value = "Hortaliza"
